$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.5762313333333333
$ws.Range("N2").Value = 1.728694
$ws.Range("O2").Value = 0.07066599920647125
$ws.Range("P2").Value = 0.07066599920647126
$ws.Range("Q2").Value = 0.2004510969242222
$ws.Range("R2").Value = 1.804059872318
$ws.Range("S2").Value = 0.07066599920647125
$ws.Range("T2").Value = 0.07066599920647126

$ws.Range("M3").Value = 3.933615666666666
$ws.Range("O3").Value = 0.4823980674067757
$ws.Range("P3").Value = 0.4823980674067757
$ws.Range("Q3").Value = 1.368369836295444
$ws.Range("S3").Value = 0.4823980674067757
$ws.Range("T3").Value = 0.4823980674067757

$ws.Range("O4").Value = 0.446935933386753
$ws.Range("P4").Value = 0.446935933386753
$ws.Range("S4").Value = 0.446935933386753
$ws.Range("T4").Value = 0.446935933386753
